# Updates cryptocurrency Price (column D) and Volume(1h) (column E) values
# on the "cryptos" worksheet, reflecting the latest scrape snapshot.
# Commit message: Updated cryptos list on Tue Sep 19 03:49:35 UTC 2023 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values look like plain numbers (e.g. "216.83", "0.0504") but must
# stay stored as text, exactly as scraped (matching trailing zeros, etc.).
# Force those specific cells to Text format before writing so Excel does not
# silently convert them into numeric values.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2: Bitcoin
$ws.Range("D2").Value = "26.859.18"
$ws.Range("E2").Value = "  +0.20%  "
# Row 3: Ethereum
$ws.Range("D3").Value = "1.636.39"
$ws.Range("E3").Value = "  -0.15%  "
# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.49%  "
# Row 5: BNB
$ws.Range("D5").Value = "216.83"
$ws.Range("E5").Value = "  -0.90%  "
# Row 6: XRP
$ws.Range("D6").Value = "0.511"
$ws.Range("E6").Value = "  +2.02%  "
# Row 7: USDC
$ws.Range("E7").Value = "  -0.41%  "
# Row 8: Cardano
$ws.Range("D8").Value = "0.255"
$ws.Range("E8").Value = "  +1.75%  "
# Row 9: Dogecoin
$ws.Range("E9").Value = "  +0.07%  "
# Row 10: Solana
$ws.Range("D10").Value = "19.92"
$ws.Range("E10").Value = "  +3.20%  "
# Row 11: TRON
$ws.Range("E11").Value = "  +0.02%  "
# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.865.42"
$ws.Range("E12").Value = "  -0.14%  "
# Row 13: WrappedEther
$ws.Range("D13").Value = "1.630.61"
$ws.Range("E13").Value = "  -0.26%  "
# Row 14: Polkadot
$ws.Range("E14").Value = "  -0.83%  "
# Row 15: Polygon
$ws.Range("E15").Value = "  +0.45%  "
# Row 16: Litecoin
$ws.Range("D16").Value = "66.67"
$ws.Range("E16").Value = "  +2.80%  "
# Row 17: WrappedBTC
$ws.Range("D17").Value = "26.860.27"
$ws.Range("E17").Value = "  +0.22%  "
# Row 18: ShibaInu
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  -0.57%  "
# Row 19: BitcoinCash
$ws.Range("D19").Value = "219.31"
$ws.Range("E19").Value = "  +1.46%  "
# Row 20: Dai
$ws.Range("E20").Value = "  -0.51%  "
# Row 21: Chainlink
$ws.Range("D21").Value = "6.73"
$ws.Range("E21").Value = "  +1.82%  "
# Row 22: Uniswap
$ws.Range("E22").Value = "  +0.68%  "
# Row 23: Toncoin
$ws.Range("D23").Value = "2.44"
$ws.Range("E23").Value = "  +3.77%  "
# Row 24: Avalanche
$ws.Range("E24").Value = "  -0.05%  "
# Row 25: Monero
$ws.Range("D25").Value = "147.05"
$ws.Range("E25").Value = "  -0.11%  "
# Row 26: BinanceUSD
$ws.Range("E26").Value = "  -0.43%  "
# Row 27: Cosmos
$ws.Range("E27").Value = "  +4.27%  "
# Row 28: Stellar
$ws.Range("E28").Value = "  +0.58%  "
# Row 29: EthereumClassic
$ws.Range("E29").Value = "  +0.22%  "
# Row 30: Hedera
$ws.Range("D30").Value = "0.0504"
$ws.Range("E30").Value = "  -0.51%  "
# Row 31: PancakeSwap
$ws.Range("E31").Value = "  -1.33%  "
# Row 32: Filecoin
$ws.Range("E32").Value = "  -1.29%  "
# Row 33: InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +0.71%  "
# Row 34: LidoDAOToken
$ws.Range("E34").Value = "  +0.77%  "
# Row 35: Maker
$ws.Range("D35").Value = "1.256.59"
$ws.Range("E35").Value = "  -0.47%  "
# Row 36: HuobiToken
$ws.Range("E36").Value = "  -0.25%  "
# Row 37: VeChain
$ws.Range("E37").Value = "  +1.92%  "
# Row 38: ImmutableX
$ws.Range("E38").Value = "  +0.65%  "
# Row 39: ARBITRUM
$ws.Range("E39").Value = "  +1.84%  "
# Row 40: PaxDollar
$ws.Range("E40").Value = "  -0.43%  "
# Row 41: TrustWalletToken
$ws.Range("D41").Value = "0.810"
$ws.Range("E41").Value = "  +0.55%  "
# Row 42: FraxShare
$ws.Range("D42").Value = "5.40"
$ws.Range("E42").Value = "  +1.19%  "
# Row 43: RocketPoolETH
$ws.Range("D43").Value = "1.779.77"
$ws.Range("E43").Value = "  +0.04%  "
# Row 44: MXToken
$ws.Range("E44").Value = "  -1.47%  "
# Row 45: Aave
$ws.Range("E45").Value = "  +2.20%  "
# Row 46: Quant
$ws.Range("D46").Value = "91.61"
$ws.Range("E46").Value = "  -0.57%  "
# Row 47: RenderToken
$ws.Range("E47").Value = "  +0.50%  "
# Row 48: BabyDogeCoin
$ws.Range("E48").Value = "  +2.73%  "
# Row 49: Cronos
$ws.Range("E49").Value = "  -0.36%  "
# Row 50: EnergySwap
$ws.Range("D50").Value = "7.64"
$ws.Range("E50").Value = "  +1.12%  "
# Row 51: Algorand
$ws.Range("D51").Value = "0.0960"
$ws.Range("E51").Value = "  -0.31%  "

Write-Host "Updated cryptos price/volume values."
